# Weekly update: a new price record (week of 2023-11-21) is reported for
# Agrícola del Norte S.A. de Arica - Albahaca. It is inserted as the new
# row 26, which pushes all the existing records (old rows 26..82) down by
# one row (new rows 27..83), growing the used range from A1:R82 to A1:R83.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 26 - shifts rows 26..82 down to 27..83.
$ws.Rows.Item(26).EntireRow.Insert()

# Fill in the new row 26 with the latest weekly record.
$ws.Cells.Item(26, 1).Value = 1
$ws.Cells.Item(26, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(26, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(26, 4).Value = 45251
$ws.Cells.Item(26, 5).Value = 15
$ws.Cells.Item(26, 6).Value = 100112052
$ws.Cells.Item(26, 7).Value = "Albahaca"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 600
$ws.Cells.Item(26, 11).Value = 800
$ws.Cells.Item(26, 12).Value = 1000
$ws.Cells.Item(26, 13).Value = 883
$ws.Cells.Item(26, 14).Value = "$/paquete"
$ws.Cells.Item(26, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(26, 16).Value = 883
$ws.Cells.Item(26, 17).Value = 1
$ws.Cells.Item(26, 18).Value = "Hortaliza"
